$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for rows 2-62.
# The value 45212 (2023-10-13) was updated to 45221 (2023-10-22).
for ($row = 2; $row -le 62; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
